$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = 111426042
$ws.Range("B6").Value = 98535
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 222498
$ws.Range("F6").Value = "Blåsippa"
$ws.Range("G6").Value = "Hepatica nobilis"
$ws.Range("H6").Value = "Schreb."
$ws.Range("Q6").Value = 550956.02874151
$ws.Range("R6").Value = 7001949.318344167
$ws.Range("Z6").Value = "15:34"
$ws.Range("AB6").Value = "15:34"

# Row 7
$ws.Range("A7").Value = 111420481
$ws.Range("B7").Value = 96381
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 219874
$ws.Range("F7").Value = "Nattviol"
$ws.Range("G7").Value = "Platanthera bifolia"
$ws.Range("H7").Value = "(L.) Rich."
$ws.Range("Q7").Value = 550727.9291679059
$ws.Range("R7").Value = 7002200.33458891
$ws.Range("Z7").Value = "11:36"
$ws.Range("AB7").Value = "11:36"

# Row 8
$ws.Range("A8").Value = 111422950
$ws.Range("B8").Value = 96348
$ws.Range("D8").Value = "VU"
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = "Knärot"
$ws.Range("G8").Value = "Goodyera repens"
$ws.Range("H8").Value = "(L.) R. Br."
$ws.Range("Q8").Value = 550849.6384025981
$ws.Range("R8").Value = 7001976.135959106
$ws.Range("Z8").Value = "13:40"
$ws.Range("AB8").Value = "13:40"

# Row 9
$ws.Range("A9").Value = 111423943
$ws.Range("B9").Value = 77267
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 6446
$ws.Range("F9").Value = "Kolflarnlav"
$ws.Range("G9").Value = "Carbonicola anthracophila"
$ws.Range("H9").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q9").Value = 550983.7409033215
$ws.Range("R9").Value = 7002058.742119073
$ws.Range("Z9").Value = "13:55"
$ws.Range("AB9").Value = "13:55"

# Row 39
$ws.Range("A39").Value = 111420869
$ws.Range("B39").Value = 89405
$ws.Range("D39").Value = "NT"
$ws.Range("E39").Value = 1202
$ws.Range("F39").Value = "Ullticka"
$ws.Range("G39").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H39").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q39").Value = 550741.937658608
$ws.Range("R39").Value = 7002115.539248759
$ws.Range("Z39").Value = "11:56"
$ws.Range("AB39").Value = "11:56"

# Row 40
$ws.Range("A40").Value = 111421476
$ws.Range("B40").Value = 77186
$ws.Range("D40").Value = "NT"
$ws.Range("E40").Value = 353
$ws.Range("F40").Value = "Dvärgbägarlav"
$ws.Range("G40").Value = "Cladonia parasitica"
$ws.Range("H40").Value = "(Hoffm.) Hoffm."
$ws.Range("Q40").Value = 550808.5392307156
$ws.Range("R40").Value = 7002084.47682756
$ws.Range("Z40").Value = "12:18"
$ws.Range("AB40").Value = "12:18"

# Row 41
$ws.Range("A41").Value = 111423255
$ws.Range("B41").Value = 98535
$ws.Range("D41").Value = "LC"
$ws.Range("E41").Value = 222498
$ws.Range("F41").Value = "Blåsippa"
$ws.Range("G41").Value = "Hepatica nobilis"
$ws.Range("H41").Value = "Schreb."
$ws.Range("Q41").Value = 550869.8149192812
$ws.Range("R41").Value = 7001960.172576382
$ws.Range("Z41").Value = "13:40"
$ws.Range("AB41").Value = "13:40"

# Row 42
$ws.Range("A42").Value = 111425789
$ws.Range("B42").Value = 99413
$ws.Range("D42").Value = "LC"
$ws.Range("E42").Value = 221235
$ws.Range("F42").Value = "Vårärt"
$ws.Range("G42").Value = "Lathyrus vernus"
$ws.Range("H42").Value = "(L.) Bernh."
$ws.Range("Q42").Value = 550964.3063330664
$ws.Range("R42").Value = 7001941.308390027
$ws.Range("Z42").Value = "15:34"
$ws.Range("AB42").Value = "15:34"

# Row 43
$ws.Range("A43").Value = 111427222
$ws.Range("B43").Value = 98535
$ws.Range("D43").Value = "LC"
$ws.Range("E43").Value = 222498
$ws.Range("F43").Value = "Blåsippa"
$ws.Range("G43").Value = "Hepatica nobilis"
$ws.Range("H43").Value = "Schreb."
$ws.Range("Q43").Value = 550980.9358707955
$ws.Range("R43").Value = 7001891.823664788
$ws.Range("Z43").Value = "16:41"
$ws.Range("AB43").Value = "16:41"

# Row 44
$ws.Range("A44").Value = 111427612
$ws.Range("B44").Value = 98446
$ws.Range("D44").Value = "LC"
$ws.Range("E44").Value = 222771
$ws.Range("F44").Value = "Svart trolldruva"
$ws.Range("G44").Value = "Actaea spicata"
$ws.Range("H44").Value = "L."
$ws.Range("Q44").Value = 550925.5549388798
$ws.Range("R44").Value = 7001928.940230627
$ws.Range("Z44").Value = "16:41"
$ws.Range("AB44").Value = "16:41"

# Row 45
$ws.Range("A45").Value = 111425793
$ws.Range("B45").Value = 98535
$ws.Range("D45").Value = "LC"
$ws.Range("E45").Value = 222498
$ws.Range("F45").Value = "Blåsippa"
$ws.Range("G45").Value = "Hepatica nobilis"
$ws.Range("H45").Value = "Schreb."
$ws.Range("Q45").Value = 550964.3063330664
$ws.Range("R45").Value = 7001941.308390027
$ws.Range("Z45").Value = "15:34"
$ws.Range("AB45").Value = "15:34"

# Row 46
$ws.Range("A46").Value = 111426315
$ws.Range("B46").Value = 89845
$ws.Range("D46").Value = "VU"
$ws.Range("E46").Value = 1209
$ws.Range("F46").Value = "Rynkskinn"
$ws.Range("G46").Value = "Phlebia centrifuga"
$ws.Range("H46").Value = "P.Karst."
$ws.Range("Q46").Value = 550958.3735980184
$ws.Range("R46").Value = 7001915.437287232
$ws.Range("Z46").Value = "16:07"
$ws.Range("AB46").Value = "16:07"

# Row 47
$ws.Range("A47").Value = 111425037
$ws.Range("B47").Value = 89845
$ws.Range("D47").Value = "VU"
$ws.Range("E47").Value = 1209
$ws.Range("F47").Value = "Rynkskinn"
$ws.Range("G47").Value = "Phlebia centrifuga"
$ws.Range("H47").Value = "P.Karst."
$ws.Range("Q47").Value = 551037.0775894802
$ws.Range("R47").Value = 7001950.143101228
$ws.Range("Z47").Value = "15:16"
$ws.Range("AB47").Value = "15:16"

# Row 67
$ws.Range("A67").Value = 111424030
$ws.Range("B67").Value = 78578
$ws.Range("D67").Value = "NT"
$ws.Range("E67").Value = 6458
$ws.Range("F67").Value = "Lunglav"
$ws.Range("G67").Value = "Lobaria pulmonaria"
$ws.Range("H67").Value = "(L.) Hoffm."
$ws.Range("Q67").Value = 551011.8102739404
$ws.Range("R67").Value = 7002059.18440557
$ws.Range("Z67").Value = "13:55"
$ws.Range("AB67").Value = "13:55"
$ws.Range("AC67").Value = "På björk"

# Row 68
$ws.Range("A68").Value = 111420370
$ws.Range("B68").Value = 96348
$ws.Range("D68").Value = "VU"
$ws.Range("E68").Value = 220787
$ws.Range("F68").Value = "Knärot"
$ws.Range("G68").Value = "Goodyera repens"
$ws.Range("H68").Value = "(L.) R. Br."
$ws.Range("Q68").Value = 550722.7296859198
$ws.Range("R68").Value = 7002214.271300747
$ws.Range("Z68").Value = "11:36"
$ws.Range("AB68").Value = "11:36"

# Row 76
$ws.Range("A76").Value = 111421391
$ws.Range("B76").Value = 78578
$ws.Range("D76").Value = "NT"
$ws.Range("E76").Value = 6458
$ws.Range("F76").Value = "Lunglav"
$ws.Range("G76").Value = "Lobaria pulmonaria"
$ws.Range("H76").Value = "(L.) Hoffm."
$ws.Range("Q76").Value = 550825.6664593286
$ws.Range("R76").Value = 7002060.778115767
$ws.Range("Z76").Value = "12:18"
$ws.Range("AB76").Value = "12:18"
$ws.Range("AC76").Value = "På asp"

# Row 77
$ws.Range("A77").Value = 111423560
$ws.Range("B77").Value = 94134
$ws.Range("D77").Value = "NT"
$ws.Range("E77").Value = 53
$ws.Range("F77").Value = "Vedtrappmossa"
$ws.Range("G77").Value = "Crossocalyx hellerianus"
$ws.Range("H77").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("Q77").Value = 550948.589889885
$ws.Range("R77").Value = 7001990.354570261
$ws.Range("Z77").Value = "13:55"
$ws.Range("AB77").Value = "13:55"

# Row 78
$ws.Range("A78").Value = 111424441
$ws.Range("B78").Value = 96348
$ws.Range("D78").Value = "VU"
$ws.Range("E78").Value = 220787
$ws.Range("F78").Value = "Knärot"
$ws.Range("G78").Value = "Goodyera repens"
$ws.Range("H78").Value = "(L.) R. Br."
$ws.Range("Q78").Value = 551054.6486243291
$ws.Range("R78").Value = 7002070.713193813
$ws.Range("Z78").Value = "00:00"
$ws.Range("AB78").Value = "00:00"

# Row 79
$ws.Range("A79").Value = 111428138
$ws.Range("B79").Value = 96348
$ws.Range("D79").Value = "VU"
$ws.Range("E79").Value = 220787
$ws.Range("F79").Value = "Knärot"
$ws.Range("G79").Value = "Goodyera repens"
$ws.Range("H79").Value = "(L.) R. Br."
$ws.Range("Q79").Value = 550809.7857848165
$ws.Range("R79").Value = 7001918.528248113
$ws.Range("Z79").Value = "17:39"
$ws.Range("AB79").Value = "17:39"

# Clear AC cells that should not be present for rows 68 and 79
$ws.Range("AC68").Value = ""
$ws.Range("AC79").Value = ""